# Daily attendance processing - 2025-12-11 15:31:38
#
# The "Recorded By" column (G) stores a comma-separated list of the
# accounts that touched a session's attendance record. For this run, the
# backend re-synced the "recorded by" history for a subset of sessions,
# which rotates the list by moving its last entry to the front (i.e. the
# most-recently-synced account is promoted to the front of the list) while
# leaving every other column untouched. Only the rows actually touched by
# that day's sync are affected; all other rows keep their original value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows (in sheet "Session Analysis Results") whose "Recorded By" value was
# resynced on this run, identified by their row number on the sheet.
$rowsToRotate = @(
    2, 3, 6, 7, 10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21, 22, 24, 26, 28,
    29, 32, 33, 36, 37, 38, 39, 40, 41, 43, 44, 45, 46, 47, 48, 50, 52, 54,
    55, 58, 59, 62, 63, 64, 65, 66, 67, 69, 70, 71, 72, 73, 74, 76, 78, 83,
    84, 85, 86, 87, 90, 92, 93, 94, 96, 99, 101, 109, 110, 111, 112, 113,
    116, 118, 119, 120, 122, 125, 127, 135, 136, 137, 138, 139, 142, 144,
    145, 146, 148, 151, 153
)

foreach ($row in $rowsToRotate) {
    $cell = $ws.Range("G$row")
    $current = [string]$cell.Value2
    $parts = $current -split ',\s*'
    if ($parts.Count -gt 1) {
        $rotated = (@($parts[-1]) + $parts[0..($parts.Count - 2)]) -join ', '
        $cell.Value2 = $rotated
    }
}
